# Update the "Sheet1" smoke-test data: bump the placeholder device-id
# suffix from 86 to 126 for the Noor.Uddin (NU) / AT rows, and flip the
# "ignore" flag from "no" to "yes" for the Noor.Uddin.* rows (5-7) now
# that their password needs to be refreshed before device creation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rename the generated Manufacturer/AuthorisedRep placeholder usernames
# (86 -> 126) for both the AT and NU rows.
$ws.Range("A3").Value = "Manufacturer126_AT"
$ws.Range("A4").Value = "AuthorisedRep126_AT"
$ws.Range("A6").Value = "Manufacturer126_NU"
$ws.Range("A7").Value = "AuthorisedRep126_NU"

# Flip "ignore" to "yes" for the Noor.Uddin Business/Manufacturer/AuthorisedRep rows
# so the smoke test updates the password and creates devices for them too.
$ws.Range("C5").Value = "yes"
$ws.Range("C6").Value = "yes"
$ws.Range("C7").Value = "yes"

# Move the active selection from A8 to C8, matching the author's last edit.
$ws.Range("C8").Select()
